# Spare_Requested_Parts.xlsx — add "Parts Code" column (CRM-2487)
#
# The export sheet has two header rows: row 1 = human-readable column
# titles, row 2 = the `{spare:...}` template placeholders substituted
# when the export runs. A new "Parts Code" column is inserted right
# before the existing "Parts Required" column (after "Age of
# Requested"), shifting every following column one slot to the right
# (old F..Q -> new G..R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, pushing "Parts Required" (and everything
# after it) one column to the right. Excel copies the left neighbour's
# column formatting automatically, same as a manual right-click ->
# Insert in the UI.
$ws.Columns("F:F").Insert()

# Match the width of the rest of that header block (the old F / new G
# column) instead of leaving the generic default width.
$ws.Columns("F:F").ColumnWidth = $ws.Columns("G:G").ColumnWidth

# Header row (row 1): plain column title, formatted like the other
# normal header cells (e.g. "Age of Requested" in E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Parts Code"

# Placeholder row (row 2): template token consumed by the export job,
# formatted like the other normal placeholder cells (e.g. E2).
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F2").Value = "{spare:part_number}"

$excel.CutCopyMode = $false
